# Yearly and Monthly Reports generation
# Applies the PPE_Report.xlsx content update:
#  - updates the accountable officer / office name in the header block
#  - adds new PPE inventory line items and two new equipment-group headers
#  - recomputes dependent Qty/Value/Depreciation numbers for the edited rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header block (row 5): accountable officer & office -------------------
$ws.Range("B5").Value = "Boris Ida Stoltenberg"
$ws.Range("C5").Value = "City Accounting Department"

# ---- Existing "Land (201) Code 1" group (row 10 header already in place) --
# Row 11: Chart holder / Sticky Notes
$ws.Range("B11").Value = "Chart holder"
$ws.Range("C11").Value = "Sticky Notes"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "2019-03-26"
$ws.Range("F11").Value = "Bailey, Joannie"
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = "I"
$ws.Range("I11").Value = 500.0
$ws.Range("J11").Value = 500.0
$ws.Range("K11").Value = 250000.0
$ws.Range("L11").Value = " "
$ws.Range("M11").Value = " "
$ws.Range("N11").Value = 450.0
$ws.Range("O11").Value = 50.0

# Row 12: Emergency light / Water Dispenser
$ws.Range("A12").Value = "01.54"
$ws.Range("B12").Value = "Emergency light"
$ws.Range("C12").Value = "Water Dispenser"
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = "2019-04-02"
$ws.Range("F12").Value = "Bailey, Joannie"
$ws.Range("G12").ClearContents()
$ws.Range("H12").Value = "I"
$ws.Range("I12").Value = 3.0
$ws.Range("J12").Value = 3.0
$ws.Range("K12").Value = 9.0
$ws.Range("L12").Value = " "
$ws.Range("M12").Value = " "
$ws.Range("N12").Value = 0.675
$ws.Range("O12").Value = 0.3

# ---- Row 13 becomes the "Office Buildings (211) Code 2" group header ------
$ws.Range("B13:O13").ClearContents()
$ws.Range("A13").Value = "Office Buildings (211) Code 2"
$ws.Range("A13").Style = $ws.Range("A10").Style

# Row 14: Chart holder / Ball Pen (new data row)
$ws.Range("A14").Value = "02.510"
$ws.Range("B14").Value = "Chart holder"
$ws.Range("C14").Value = "Ball Pen"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "2019-03-26"
$ws.Range("F14").Value = "Bailey, Joannie"
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = "I"
$ws.Range("I14").Value = 200.0
$ws.Range("J14").Value = 200.0
$ws.Range("K14").Value = 40000.0
$ws.Range("L14").Value = " "
$ws.Range("M14").Value = " "
$ws.Range("N14").Value = 180.0
$ws.Range("O14").Value = 20.0

# ---- Row 15: new "Hospitals and Health Centers (213) Code 4" group header -
$ws.Range("A15").Value = "Hospitals and Health Centers (213) Code 4"
$ws.Range("A15").Style = $ws.Range("A10").Style

# Row 16: Instructional Material / Water bottles (new data row)
$ws.Range("A16").Value = "04.42"
$ws.Range("B16").Value = "Instructional Material"
$ws.Range("C16").Value = "Water bottles"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "2019-04-02"
$ws.Range("F16").Value = "Bailey, Joannie"
$ws.Range("H16").Value = "I"
$ws.Range("I16").Value = 69.0
$ws.Range("J16").Value = 69.0
$ws.Range("K16").Value = 4761.0
$ws.Range("L16").Value = " "
$ws.Range("M16").Value = " "
$ws.Range("N16").Value = 31.05
$ws.Range("O16").Value = 6.9

# ---- Column widths: re-fit to the new (longer / shorter) content ----------
$ws.Columns("A:O").AutoFit()
